$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TonKho")

# Update QUANTITY column (E) values for the inventory refresh (29/12/2023).
# These cells store their numeric-looking quantities as TEXT (shared strings),
# matching the rest of the sheet, so plain Range.Value assignment (which
# auto-coerces a numeric-looking string to a real number) is avoided:
#  - E5 / E7 are repointed to existing text cells that already hold the
#    target text ("10" and "5" respectively) via Copy/Paste, which keeps
#    the text type intact.
#  - E6 needs a brand-new text value ("39"), so a helper cell builds it
#    with a TEXT() formula (always text) and that result is pasted as a
#    value into E6, then the helper is cleared.

$ws.Range("E2").Copy($ws.Range("E5"))
$ws.Range("E4").Copy($ws.Range("E7"))

$helper = $ws.Range("Z1")
$helper.Formula = '=TEXT(39,"0")'
$helper.Copy()
$ws.Range("E6").PasteSpecial(-4163)
$helper.ClearContents()

$excel.CutCopyMode = 0
